$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-07-07"

# Update header label cell (shared string "2022 (through 07-06)" -> "2022 (through 07-07)")
$ws.Range("I1").Value = "2022 (through 07-07)"

# Update July (row 8) Total column value
$ws.Range("I8").Value = 34

# Update yearly Total row (row 14) Total column value
$ws.Range("I14").Value = 840
